# Auto-generated: update market-price / profit columns (H, I, J, K, L, M, N)
# on multiple item leve rows across sheets, per the scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 3141
$ws.Range("I7").Value = 282
$ws.Range("K7").Value = 282
$ws.Range("M7").Value = -170

$ws.Range("H14").Value = 3141
$ws.Range("I14").Value = 282
$ws.Range("K14").Value = 282
$ws.Range("M14").Value = -91

$ws.Range("H19").Value = 804.8261
$ws.Range("I19").Value = 462.2857
$ws.Range("J19").Value = 954.6875
$ws.Range("K19").Value = 462.2857
$ws.Range("L19").Value = 954.6875
$ws.Range("M19").Value = -287.2857
$ws.Range("N19").Value = -1304.6875

$ws.Range("H64").Value = 3484.3125
$ws.Range("I64").Value = 3354.125
$ws.Range("J64").Value = 3874.875
$ws.Range("K64").Value = 3354.125
$ws.Range("L64").Value = 3874.875
$ws.Range("M64").Value = -3106.125
$ws.Range("N64").Value = -4370.875

$ws.Range("H67").Value = 3484.3125
$ws.Range("I67").Value = 3354.125
$ws.Range("J67").Value = 3874.875
$ws.Range("K67").Value = 3354.125
$ws.Range("L67").Value = 3874.875
$ws.Range("M67").Value = -2496.125
$ws.Range("N67").Value = -5590.875

$ws.Range("H76").Value = 3501.1304
$ws.Range("I76").Value = 3386.05
$ws.Range("K76").Value = 3386.05
$ws.Range("M76").Value = -3071.05

$ws.Range("H79").Value = 3501.1304
$ws.Range("I79").Value = 3386.05
$ws.Range("K79").Value = 3386.05
$ws.Range("M79").Value = -2294.05

$ws.Range("H100").Value = 1946.909
$ws.Range("I100").Value = 1065.3334
$ws.Range("K100").Value = 1065.3334
$ws.Range("M100").Value = -524.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7047.519
$ws.Range("I32").Value = 5426.625
$ws.Range("J32").Value = 13963.333
$ws.Range("K32").Value = 5426.625
$ws.Range("L32").Value = 13963.333
$ws.Range("M32").Value = -5139.625
$ws.Range("N32").Value = -14537.333

$ws.Range("H63").Value = 4715.7144
$ws.Range("I63").Value = 2501.6667
$ws.Range("J63").Value = 18000
$ws.Range("K63").Value = 2501.6667
$ws.Range("L63").Value = 18000
$ws.Range("M63").Value = -1815.6667
$ws.Range("N63").Value = -19372

$ws.Range("H66").Value = 4715.7144
$ws.Range("I66").Value = 2501.6667
$ws.Range("J66").Value = 18000
$ws.Range("K66").Value = 12508.3335
$ws.Range("L66").Value = 90000
$ws.Range("M66").Value = -9076.333500000001
$ws.Range("N66").Value = -96864

$ws.Range("H88").Value = 5510.5713
$ws.Range("I88").Value = 13976.5
$ws.Range("J88").Value = 2124.2
$ws.Range("K88").Value = 13976.5
$ws.Range("L88").Value = 2124.2
$ws.Range("M88").Value = -13570.5
$ws.Range("N88").Value = -2936.2

$ws.Range("H91").Value = 5510.5713
$ws.Range("I91").Value = 13976.5
$ws.Range("J91").Value = 2124.2
$ws.Range("K91").Value = 13976.5
$ws.Range("L91").Value = 2124.2
$ws.Range("M91").Value = -12572.5
$ws.Range("N91").Value = -4932.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1072
$ws.Range("I20").Value = 1068.2307
$ws.Range("J20").Value = 1081.8
$ws.Range("K20").Value = 1068.2307
$ws.Range("L20").Value = 1081.8
$ws.Range("M20").Value = -821.2307000000001
$ws.Range("N20").Value = -1575.8

$ws.Range("H105").Value = 5293.5
$ws.Range("I105").Value = 4948.5
$ws.Range("K105").Value = 4948.5
$ws.Range("M105").Value = -3201.5

$ws.Range("H107").Value = 2095.889
$ws.Range("I107").Value = 2069.2307
$ws.Range("J107").Value = 2165.2
$ws.Range("K107").Value = 2069.2307
$ws.Range("L107").Value = 2165.2
$ws.Range("M107").Value = -149.2307000000001
$ws.Range("N107").Value = -6005.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3951.111
$ws.Range("I62").Value = 3572.1428
$ws.Range("K62").Value = 3572.1428
$ws.Range("M62").Value = -2948.1428

$ws.Range("H65").Value = 3951.111
$ws.Range("I65").Value = 3572.1428
$ws.Range("K65").Value = 17860.714
$ws.Range("M65").Value = -14740.714

$ws.Range("H86").Value = 1526.2
$ws.Range("I86").Value = 1166.6666
$ws.Range("J86").Value = 1820.3636
$ws.Range("K86").Value = 1166.6666
$ws.Range("L86").Value = 1820.3636
$ws.Range("M86").Value = -43.66660000000002
$ws.Range("N86").Value = -4066.3636

$ws.Range("H89").Value = 1526.2
$ws.Range("I89").Value = 1166.6666
$ws.Range("J89").Value = 1820.3636
$ws.Range("K89").Value = 5833.333000000001
$ws.Range("L89").Value = 9101.817999999999
$ws.Range("M89").Value = -217.3330000000005
$ws.Range("N89").Value = -20333.818

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 263.6
$ws.Range("I6").Value = 39.333332
$ws.Range("J6").Value = 600
$ws.Range("K6").Value = 117.999996
$ws.Range("L6").Value = 1800
$ws.Range("M6").Value = -4.999995999999996
$ws.Range("N6").Value = -2026

$ws.Range("H131").Value = 675.2
$ws.Range("I131").Value = 372.85715
$ws.Range("J131").Value = 755.56964
$ws.Range("K131").Value = 1118.57145
$ws.Range("L131").Value = 2266.70892
$ws.Range("M131").Value = 3921.42855
$ws.Range("N131").Value = -12346.70892

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5741.8726
$ws.Range("I70").Value = 5591.7
$ws.Range("J70").Value = 5853.1113
$ws.Range("K70").Value = 5591.7
$ws.Range("L70").Value = 5853.1113
$ws.Range("M70").Value = -5321.7
$ws.Range("N70").Value = -6393.1113

$ws.Range("H73").Value = 5741.8726
$ws.Range("I73").Value = 5591.7
$ws.Range("J73").Value = 5853.1113
$ws.Range("K73").Value = 5591.7
$ws.Range("L73").Value = 5853.1113
$ws.Range("M73").Value = -4655.7
$ws.Range("N73").Value = -7725.1113

$ws.Range("H80").Value = 7350
$ws.Range("I80").Value = 18465
$ws.Range("J80").Value = 3645
$ws.Range("K80").Value = 18465
$ws.Range("L80").Value = 3645
$ws.Range("M80").Value = -17467
$ws.Range("N80").Value = -5641

$ws.Range("H83").Value = 7350
$ws.Range("I83").Value = 18465
$ws.Range("J83").Value = 3645
$ws.Range("K83").Value = 92325
$ws.Range("L83").Value = 18225
$ws.Range("M83").Value = -87333
$ws.Range("N83").Value = -28209

$ws.Range("H122").Value = 7520.8887
$ws.Range("I122").Value = 11740
$ws.Range("J122").Value = 2247
$ws.Range("K122").Value = 35220
$ws.Range("L122").Value = 6741
$ws.Range("M122").Value = -32770
$ws.Range("N122").Value = -11641

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5609.2856
$ws.Range("I7").Value = 3820.5
$ws.Range("J7").Value = 7235.4546
$ws.Range("K7").Value = 3820.5
$ws.Range("L7").Value = 7235.4546
$ws.Range("M7").Value = -3708.5
$ws.Range("N7").Value = -7459.4546

$ws.Range("H16").Value = 963.0476
$ws.Range("I16").Value = 998.3125
$ws.Range("J16").Value = 850.2
$ws.Range("K16").Value = 998.3125
$ws.Range("L16").Value = 850.2
$ws.Range("M16").Value = -828.3125
$ws.Range("N16").Value = -1190.2

$ws.Range("H20").Value = 5340.8237
$ws.Range("I20").Value = 1000
$ws.Range("J20").Value = 5919.6
$ws.Range("K20").Value = 1000
$ws.Range("L20").Value = 5919.6
$ws.Range("M20").Value = -774
$ws.Range("N20").Value = -6371.6

$ws.Range("H40").Value = 4031.3928
$ws.Range("I40").Value = 3722.9048
$ws.Range("J40").Value = 4956.857
$ws.Range("K40").Value = 3722.9048
$ws.Range("L40").Value = 4956.857
$ws.Range("M40").Value = -3586.9048
$ws.Range("N40").Value = -5228.857

$ws.Range("H126").Value = 5609.2856
$ws.Range("I126").Value = 3820.5
$ws.Range("J126").Value = 7235.4546
$ws.Range("K126").Value = 11461.5
$ws.Range("L126").Value = 21706.3638
$ws.Range("M126").Value = -8991.5
$ws.Range("N126").Value = -26646.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 2971.5715
$ws.Range("J4").Value = 2971.5715
$ws.Range("L4").Value = 2971.5715
$ws.Range("N4").Value = -3197.5715
